$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number by
# Excel's smart entry -- force text format first, write value, then drop the
# format override again so the cell keeps its original (default/general) style.

$textCells = @("D5", "D6", "D7", "D9", "D10", "D12", "D14", "D15", "D17", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D48", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.067.72"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "2.048.81"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "248.97"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").Value = "0.669"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "59.16"
$ws.Range("E7").Value = "  +7.91%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.387"
$ws.Range("E9").Value = "  +0.84%  "

$ws.Range("D10").Value = "0.0790"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("D12").Value = "15.94"
$ws.Range("E12").Value = "  +5.85%  "

$ws.Range("D13").Value = "2.352.00"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Value = "0.832"
$ws.Range("E14").Value = "  +1.78%  "

$ws.Range("D15").Value = "5.71"
$ws.Range("E15").Value = "  +6.76%  "

$ws.Range("D16").Value = "2.049.71"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").Value = "18.64"
$ws.Range("E17").Value = "  +30.82%  "

$ws.Range("D18").Value = "37.054.52"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").Value = "75.39"
$ws.Range("E19").Value = "  +2.56%  "

$ws.Range("D20").Value = "0.0₃0901"

$ws.Range("D21").Value = "5.41"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "237.81"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  +10.54%  "

$ws.Range("D26").Value = "9.55"
$ws.Range("E26").Value = "  +5.71%  "

$ws.Range("D27").Value = "168.62"
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("D28").Value = "20.08"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("E29").Value = "  +0.56%  "

$ws.Range("D30").Value = "1.12"
$ws.Range("E30").Value = "  +6.24%  "

$ws.Range("D31").Value = "4.78"
$ws.Range("E31").Value = "  +3.78%  "

$ws.Range("D32").Value = "0.0626"
$ws.Range("E32").Value = "  -0.52%  "

$ws.Range("D33").Value = "4.51"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("D34").Value = "0.0892"
$ws.Range("E34").Value = "  -0.41%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  -3.28%  "

$ws.Range("D37").Value = "1.74"
$ws.Range("E37").Value = "  -1.86%  "

$ws.Range("E38").Value = "  +4.31%  "

$ws.Range("D39").Value = "1.33"
$ws.Range("E39").Value = "  -1.45%  "

$ws.Range("D40").Value = "3.10"
$ws.Range("E40").Value = "  +10.68%  "

$ws.Range("D41").Value = "5.08"
$ws.Range("E41").Value = "  +23.34%  "

$ws.Range("D42").Value = "17.66"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").Value = "0.0223"
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("D45").Value = "96.77"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("E46").Value = "  +4.16%  "

$ws.Range("D47").Value = "1.288.85"
$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("D48").Value = "3.81"
$ws.Range("E48").Value = "  -8.75%  "

$ws.Range("D49").Value = "2.87"
$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("D50").Value = "6.81"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").Value = "2.234.51"
$ws.Range("E51").Value = "  -0.90%  "

# Restore the default style on the forced-text cells (keeps the saved file free
# of any lingering custom cell style now that the text is locked in).
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
